$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Rows 27-31 (columns A and B) get shifted: the "Microsoft"/"Microsoft" row
# currently at row 27 moves down to row 31, and rows 28-31 shift up to 27-30.
$ws.Range("A27").Value = "Aspire"
$ws.Range("B27").Value = "Acer"

$ws.Range("A28").Value = "ConceptD"
$ws.Range("B28").Value = "Acer"

$ws.Range("A29").Value = "Dynabook"
$ws.Range("B29").Value = "Dynabook Toshiba"

$ws.Range("A30").Value = "HP"
$ws.Range("B30").Value = "HP"

$ws.Range("A31").Value = "Microsoft"
$ws.Range("B31").Value = "Microsoft"

$wb.Save()
